# Saldo_guide.xlsx update script
# - Rename sheet tab (reflects new extraction date 2024-10-31 vs 2024-10-30)
# - Bump every "Dt. Referencia" (column G) value from 45595 to 45596 (one day later)
# - Apply the handful of balance corrections for specific accounts (columns D/E/H)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab to match the new reference date
$ws.Name = "IClientBalance-20241031-092739-"

# Shift the reference date (column G) for every data row (2..274) by one day
for ($r = 2; $r -le 274; $r++) {
    $ws.Cells.Item($r, 7).Value = 45596
}

# Row 107 (F=5534961660): Saldo Previsto / Vl. Total corrected
$ws.Cells.Item(107, 5).Value = 66347.45
$ws.Cells.Item(107, 8).Value = 66347.45

# Row 108 (F=15107385761): Vl. Projetado / Vl. Total corrected
$ws.Cells.Item(108, 4).Value = -161.34
$ws.Cells.Item(108, 8).Value = -161.34

# Row 114 (F=54700426772): Saldo Previsto / Vl. Total corrected
$ws.Cells.Item(114, 5).Value = 27503.59
$ws.Cells.Item(114, 8).Value = 27503.59

# Row 115 (F=2705342109): Saldo Previsto / Vl. Total corrected
$ws.Cells.Item(115, 5).Value = 0
$ws.Cells.Item(115, 8).Value = 0

# Row 143 (F=6733154447): Saldo Previsto / Vl. Total corrected
$ws.Cells.Item(143, 5).Value = 123768.59
$ws.Cells.Item(143, 8).Value = 123768.59

# Row 165 (F=92296580149): Vl. Projetado corrected, Vl. Total recomputed
$ws.Cells.Item(165, 4).Value = -776.46
$ws.Cells.Item(165, 8).Value = -672.81

# Row 189 (F=82395489972): Vl. Projetado corrected, Vl. Total recomputed
$ws.Cells.Item(189, 4).Value = -43.79
$ws.Cells.Item(189, 8).Value = 0

# Row 224 (F=27734048668): Vl. Projetado / Saldo Previsto corrected
$ws.Cells.Item(224, 4).Value = 0
$ws.Cells.Item(224, 5).Value = 10559.36
